# "added data to the excel" -- populate Sheet1 with test-account rows
# (email / password pairs) and wire up hyperlinks on the three cells
# that originally carried them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 1 ---------------------------------------------------------
$ws.Range("A1").Value = "selauto1@test.com "
$ws.Range("B1").Value = "Pass@123"

# --- row 2 ---------------------------------------------------------
$ws.Range("A2").Value = "selauto2@test.com "
$ws.Range("B2").Value = "Pass@123"

# --- row 3 ---------------------------------------------------------
$ws.Range("A3").Value = "testautouser@gmail.com"
$ws.Range("B3").Value = "Pass@123"

# --- row 4 ---------------------------------------------------------
$ws.Range("A4").Value = "testautouser2@gmail.com"
$ws.Range("B4").Value = "Pass@123"

# --- row 5 ---------------------------------------------------------
$ws.Range("A5").Value = "dummy@test123"
$ws.Range("B5").Value = "Wrong PW"

# --- hyperlinks (also applies the built-in "Hyperlink" style) ------
[void]$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:testautouser@gmail.com")
[void]$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:testautouser2@gmail.com")
[void]$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:dummy@test123")

# --- column A width (matches authored width as closely as the host
#     rounding allows) ------------------------------------------------
$ws.Columns("A").ColumnWidth = 26

# --- restore the selection left on the sheet -----------------------
[void]$ws.Range("B8").Select()
